# Insert two new rows at 261, pushing the existing rows 261..356 down to 263..358.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("261:262").Insert()

# --- New row 261 ---
$ws.Cells.Item(261, 1).Value = 10
$ws.Cells.Item(261, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(261, 3).Value = "La Araucanía"
$ws.Cells.Item(261, 4).Value = 44988
$ws.Cells.Item(261, 5).Value = 9
$ws.Cells.Item(261, 6).Value = 100112052
$ws.Cells.Item(261, 7).Value = "Albahaca"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 50
$ws.Cells.Item(261, 11).Value = 5000
$ws.Cells.Item(261, 12).Value = 6000
$ws.Cells.Item(261, 13).Value = 5600
$ws.Cells.Item(261, 14).Value = "$/paquete"
$ws.Cells.Item(261, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(261, 16).Value = 5600
$ws.Cells.Item(261, 17).Value = 1
$ws.Cells.Item(261, 18).Value = "Hortaliza"

# --- New row 262 ---
$ws.Cells.Item(262, 1).Value = 10
$ws.Cells.Item(262, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(262, 3).Value = "La Araucanía"
$ws.Cells.Item(262, 4).Value = 44988
$ws.Cells.Item(262, 5).Value = 9
$ws.Cells.Item(262, 6).Value = 100112052
$ws.Cells.Item(262, 7).Value = "Albahaca"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 40
$ws.Cells.Item(262, 11).Value = 5000
$ws.Cells.Item(262, 12).Value = 6000
$ws.Cells.Item(262, 13).Value = 5500
$ws.Cells.Item(262, 14).Value = "$/paquete"
$ws.Cells.Item(262, 15).Value = "Región del Maule"
$ws.Cells.Item(262, 16).Value = 5500
$ws.Cells.Item(262, 17).Value = 1
$ws.Cells.Item(262, 18).Value = "Hortaliza"

Write-Output "done"
